$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2235649546827795
$ws.Range("C2").Value = 0.5256797583081571
$ws.Range("J2").Value = 0.00906344410876133
$ws.Range("P2").Value = 0.1812688821752266
$ws.Range("S2").Value = 0.06042296072507553
$ws.Range("B3").Value = 0.01666666666666667
$ws.Range("C3").Value = 0.02777777777777778
$ws.Range("J3").Value = 0.02222222222222222
$ws.Range("P3").Value = 0.7666666666666667
$ws.Range("S3").Value = 0.1666666666666667
$ws.Range("J4").Value = 0.0392156862745098
$ws.Range("P4").Value = 0.7647058823529411
$ws.Range("S4").Value = 0.196078431372549
$ws.Range("B6").Value = 0.03669724770642202
$ws.Range("D6").Value = 0.004587155963302753
$ws.Range("F6").Value = 0.04587155963302753
$ws.Range("J6").Value = 0.1834862385321101
$ws.Range("O6").Value = 0.02752293577981652
$ws.Range("Q6").Value = 0.2339449541284404
$ws.Range("R6").Value = 0.0963302752293578
$ws.Range("S6").Value = 0.3715596330275229
$ws.Range("B7").Value = 0.1421319796954315
$ws.Range("D7").Value = 0.005076142131979695
$ws.Range("E7").Value = 0.005076142131979695
$ws.Range("F7").Value = 0.07614213197969544
$ws.Range("J7").Value = 0.1015228426395939
$ws.Range("O7").Value = 0.01015228426395939
$ws.Range("Q7").Value = 0.2131979695431472
$ws.Range("R7").Value = 0.09644670050761421
$ws.Range("S7").Value = 0.350253807106599
$ws.Range("B8").Value = 0.0915032679738562
$ws.Range("D8").Value = 0.02178649237472767
$ws.Range("F8").Value = 0.04357298474945534
$ws.Range("J8").Value = 0.1111111111111111
$ws.Range("O8").Value = 0.01525054466230937
$ws.Range("Q8").Value = 0.2004357298474946
$ws.Range("R8").Value = 0.08496732026143791
$ws.Range("S8").Value = 0.4313725490196079
$ws.Range("B9").Value = 0.14
$ws.Range("D9").Value = 0.01
$ws.Range("F9").Value = 0.065
$ws.Range("J9").Value = 0.11
$ws.Range("O9").Value = 0.01
$ws.Range("Q9").Value = 0.17
$ws.Range("R9").Value = 0.08
$ws.Range("S9").Value = 0.415
$ws.Range("B10").Value = 0.1125879593432369
$ws.Range("D10").Value = 0.0289288506645817
$ws.Range("F10").Value = 0.08444096950742767
$ws.Range("J10").Value = 0.1321344800625489
$ws.Range("O10").Value = 0.01016419077404222
$ws.Range("Q10").Value = 0.2040656763096169
$ws.Range("R10").Value = 0.08444096950742767
$ws.Range("S10").Value = 0.343236903831118
$ws.Range("G11").Value = 0.1457627118644068
$ws.Range("J11").Value = 0.07796610169491526
$ws.Range("K11").Value = 0.1966101694915254
$ws.Range("L11").Value = 0.5423728813559322
$ws.Range("S11").Value = 0.03728813559322034
$ws.Range("G12").Value = 0.7453416149068323
$ws.Range("J12").Value = 0.1925465838509317
$ws.Range("L12").Value = 0.01863354037267081
$ws.Range("S12").Value = 0.04347826086956522
$ws.Range("G13").Value = 0.6851851851851852
$ws.Range("J13").Value = 0.2777777777777778
$ws.Range("S13").Value = 0.03703703703703703
$ws.Range("S14").Value = 1
$ws.Range("F15").Value = 0.01587301587301587
$ws.Range("H15").Value = 0.1481481481481481
$ws.Range("I15").Value = 0.06349206349206349
$ws.Range("J15").Value = 0.3915343915343915
$ws.Range("K15").Value = 0.07407407407407407
$ws.Range("M15").Value = 0.01587301587301587
$ws.Range("O15").Value = 0.02116402116402116
$ws.Range("S15").Value = 0.2698412698412698
$ws.Range("F16").Value = 0.004464285714285714
$ws.Range("H16").Value = 0.1651785714285714
$ws.Range("I16").Value = 0.1160714285714286
$ws.Range("J16").Value = 0.3928571428571428
$ws.Range("K16").Value = 0.1160714285714286
$ws.Range("M16").Value = 0.01785714285714286
$ws.Range("O16").Value = 0.05357142857142857
$ws.Range("S16").Value = 0.1339285714285714
$ws.Range("F17").Value = 0.01257861635220126
$ws.Range("H17").Value = 0.1886792452830189
$ws.Range("I17").Value = 0.1048218029350105
$ws.Range("J17").Value = 0.4150943396226415
$ws.Range("K17").Value = 0.06708595387840671
$ws.Range("M17").Value = 0.02306079664570231
$ws.Range("O17").Value = 0.04612159329140461
$ws.Range("S17").Value = 0.1425576519916142
$ws.Range("F18").Value = 0.009852216748768473
$ws.Range("H18").Value = 0.2216748768472906
$ws.Range("I18").Value = 0.09852216748768473
$ws.Range("J18").Value = 0.3448275862068966
$ws.Range("K18").Value = 0.07881773399014778
$ws.Range("M18").Value = 0.03448275862068965
$ws.Range("O18").Value = 0.06896551724137931
$ws.Range("S18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.007009345794392523
$ws.Range("H19").Value = 0.205607476635514
$ws.Range("I19").Value = 0.07476635514018691
$ws.Range("J19").Value = 0.3722741433021807
$ws.Range("K19").Value = 0.1121495327102804
$ws.Range("M19").Value = 0.02258566978193146
$ws.Range("N19").Value = 0.000778816199376947
$ws.Range("O19").Value = 0.06542056074766354
$ws.Range("S19").Value = 0.1394080996884735

Write-Output "Applied 108 cell updates to team specific time matrix"
